# Update market-price derived columns (H-N) on several leve rows across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled-runner
# refresh of Chocobo_Profits market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H129").Value = 833.24
$ws.Range("J129").Value = 899.9773
$ws.Range("L129").Value = 2699.9319
$ws.Range("N129").Value = -12699.9319

$ws.Range("H137").Value = 3710.4614
$ws.Range("I137").Value = 2694.5
$ws.Range("J137").Value = 5336
$ws.Range("K137").Value = 8083.5
$ws.Range("L137").Value = 16008
$ws.Range("M137").Value = -5533.5
$ws.Range("N137").Value = -21108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4374.949
$ws.Range("I32").Value = 3750.8628
$ws.Range("J32").Value = 8353.5
$ws.Range("K32").Value = 3750.8628
$ws.Range("L32").Value = 8353.5
$ws.Range("M32").Value = -3463.8628
$ws.Range("N32").Value = -8927.5

$ws.Range("H61").Value = 1338
$ws.Range("I61").Value = 1386.3
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 1386.3
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -1174.3
$ws.Range("N61").Value = -1624

$ws.Range("H74").Value = 7923.2856
$ws.Range("I74").Value = 37700
$ws.Range("J74").Value = 2960.5
$ws.Range("K74").Value = 37700
$ws.Range("L74").Value = 2960.5
$ws.Range("M74").Value = -36826
$ws.Range("N74").Value = -4708.5

$ws.Range("H77").Value = 7923.2856
$ws.Range("I77").Value = 37700
$ws.Range("J77").Value = 2960.5
$ws.Range("K77").Value = 188500
$ws.Range("L77").Value = 14802.5
$ws.Range("M77").Value = -184132
$ws.Range("N77").Value = -23538.5

$ws.Range("H102").Value = 2389.2
$ws.Range("I102").Value = 1968.4
$ws.Range("K102").Value = 1968.4
$ws.Range("M102").Value = -346.4000000000001

$ws.Range("H110").Value = 1483.6666
$ws.Range("I110").Value = 1481.1538
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 1481.1538
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 563.8462
$ws.Range("N110").Value = -5590

$ws.Range("H132").Value = 3153.261
$ws.Range("I132").Value = 1972.7273
$ws.Range("K132").Value = 5918.1819
$ws.Range("M132").Value = -3388.1819

$ws.Range("H136").Value = 1338
$ws.Range("I136").Value = 1386.3
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 4158.9
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -1608.9
$ws.Range("N136").Value = -8700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5401363.5
$ws.Range("I7").Value = 15002500
$ws.Range("J7").Value = 3655702.2
$ws.Range("K7").Value = 15002500
$ws.Range("L7").Value = 3655702.2
$ws.Range("M7").Value = -15002387
$ws.Range("N7").Value = -3655928.2

$ws.Range("H103").Value = 38000
$ws.Range("J103").Value = 38000
$ws.Range("L103").Value = 38000
$ws.Range("N103").Value = -40344

$ws.Range("H105").Value = 1600.2273
$ws.Range("J105").Value = 1619.25
$ws.Range("L105").Value = 1619.25
$ws.Range("N105").Value = -5113.25

$ws.Range("H107").Value = 2211.1875
$ws.Range("I107").Value = 2098.5
$ws.Range("K107").Value = 2098.5
$ws.Range("M107").Value = -178.5

$ws.Range("H134").Value = 2236.55
$ws.Range("I134").Value = 1534.4839
$ws.Range("J134").Value = 4654.778
$ws.Range("K134").Value = 4603.4517
$ws.Range("L134").Value = 13964.334
$ws.Range("M134").Value = -2068.4517
$ws.Range("N134").Value = -19034.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 9766.333000000001
$ws.Range("I12").Value = 2150
$ws.Range("J12").Value = 24999
$ws.Range("K12").Value = 2150
$ws.Range("L12").Value = 24999
$ws.Range("M12").Value = -1980
$ws.Range("N12").Value = -25339

$ws.Range("H64").Value = 37270.75
$ws.Range("J64").Value = 37270.75
$ws.Range("L64").Value = 37270.75
$ws.Range("N64").Value = -37766.75

$ws.Range("H67").Value = 37270.75
$ws.Range("J67").Value = 37270.75
$ws.Range("L67").Value = 37270.75
$ws.Range("N67").Value = -38986.75

$ws.Range("H132").Value = 2430.1292
$ws.Range("I132").Value = 1795.619
$ws.Range("J132").Value = 3762.6
$ws.Range("K132").Value = 5386.857
$ws.Range("L132").Value = 11287.8
$ws.Range("M132").Value = -2856.857
$ws.Range("N132").Value = -16347.8

$ws.Range("H134").Value = 8161.1177
$ws.Range("I134").Value = 9561.583000000001
$ws.Range("J134").Value = 4800
$ws.Range("K134").Value = 28684.749
$ws.Range("L134").Value = 14400
$ws.Range("M134").Value = -26149.749
$ws.Range("N134").Value = -19470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1515.3846
$ws.Range("J17").Value = 1833.3334
$ws.Range("L17").Value = 5500.0002
$ws.Range("N17").Value = -5838.0002

$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

$ws.Range("H57").Value = 2755.2222
$ws.Range("I57").Value = 600
$ws.Range("J57").Value = 3024.625
$ws.Range("K57").Value = 1800
$ws.Range("L57").Value = 9073.875
$ws.Range("M57").Value = -1241
$ws.Range("N57").Value = -10191.875

$ws.Range("H113").Value = 660.1739
$ws.Range("I113").Value = 653.5454999999999
$ws.Range("K113").Value = 1960.6365
$ws.Range("M113").Value = 209.3635000000002

$ws.Range("H129").Value = 3770.1
$ws.Range("I129").Value = 3624
$ws.Range("J129").Value = 4111
$ws.Range("K129").Value = 10872
$ws.Range("L129").Value = 12333
$ws.Range("M129").Value = -5872
$ws.Range("N129").Value = -22333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3310.5334
$ws.Range("I132").Value = 2055
$ws.Range("K132").Value = 6165
$ws.Range("M132").Value = -3635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1798.2693
$ws.Range("I82").Value = 662.75
$ws.Range("J82").Value = 2771.5715
$ws.Range("K82").Value = 662.75
$ws.Range("L82").Value = 2771.5715
$ws.Range("M82").Value = -301.75
$ws.Range("N82").Value = -3493.5715

$ws.Range("H85").Value = 1798.2693
$ws.Range("I85").Value = 662.75
$ws.Range("J85").Value = 2771.5715
$ws.Range("K85").Value = 662.75
$ws.Range("L85").Value = 2771.5715
$ws.Range("M85").Value = 585.25
$ws.Range("N85").Value = -5267.5715

$ws.Range("H132").Value = 4584.387
$ws.Range("I132").Value = 1627.2632
$ws.Range("K132").Value = 4881.7896
$ws.Range("M132").Value = -2351.7896

$ws.Range("H136").Value = 3382.7878
$ws.Range("I136").Value = 1104.8235
$ws.Range("J136").Value = 5803.125
$ws.Range("K136").Value = 3314.4705
$ws.Range("L136").Value = 17409.375
$ws.Range("M136").Value = -764.4704999999999
$ws.Range("N136").Value = -22509.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 16266.667
$ws.Range("J74").Value = 16266.667
$ws.Range("L74").Value = 16266.667
$ws.Range("N74").Value = -18138.667

$ws.Range("H77").Value = 16266.667
$ws.Range("J77").Value = 16266.667
$ws.Range("L77").Value = 48800.001
$ws.Range("N77").Value = -58160.001

$ws.Range("H81").Value = 1975
$ws.Range("I81").Value = 1975
$ws.Range("K81").Value = 3950
$ws.Range("M81").Value = -2889

$ws.Range("H84").Value = 1975
$ws.Range("I84").Value = 1975
$ws.Range("K84").Value = 19750
$ws.Range("M84").Value = -14446

$ws.Range("H109").Value = 27776.5
$ws.Range("J109").Value = 27776.5
$ws.Range("L109").Value = 27776.5
$ws.Range("N109").Value = -30550.5

$ws.Range("H136").Value = 4260.5386
$ws.Range("I136").Value = 2433.65
$ws.Range("K136").Value = 7300.950000000001
$ws.Range("M136").Value = -4750.950000000001
